$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Language Key value on row 6 (Global vendor request) from "EN" to "EN,English".
#    Set this first so the new shared string "EN,English" is appended before the
#    other newly-introduced strings (matches the order new strings were authored).
$ws.Range("R6").Value = "EN,English"

# 2. Add the new "Company Trading Partner" column (BM) with header + "AE01" values
#    for each of the five data rows.
$ws.Range("BM1").Value = "Company Trading Partner"
$ws.Range("BM2").Value = "AE01"
$ws.Range("BM3").Value = "AE01"
$ws.Range("BM4").Value = "AE01"
$ws.Range("BM5").Value = "AE01"
$ws.Range("BM6").Value = "AE01"

# 3. Match formatting of the new column to the existing table: data rows inherit the
#    formatting already used by the adjacent "Bank Bearer" column (BL).
$ws.Range("BL2:BL6").Copy() | Out-Null
$ws.Range("BM2:BM6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4. Leave the selection where the author left it after entering the new column data.
$ws.Range("BM3:BM6").Select() | Out-Null
